$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Row 4 (year header 2022) - same style as K4/L4
$ws.Range("L4").Copy()
$ws.Range("M4").PasteSpecial($xlPasteFormats)
$ws.Range("M4").Value = 2022

$ws.Range("K5").Copy()
$ws.Range("M5").PasteSpecial($xlPasteFormats)
$ws.Range("M5").Value = 24.6

$ws.Range("K6").Copy()
$ws.Range("M6").PasteSpecial($xlPasteFormats)
$ws.Range("M6").Value = 40.7

$ws.Range("K7").Copy()
$ws.Range("M7").PasteSpecial($xlPasteFormats)
$ws.Range("M7").Value = 20.7

$ws.Range("K8").Copy()
$ws.Range("M8").PasteSpecial($xlPasteFormats)
$ws.Range("M8").Value = 26.6

$ws.Range("K9").Copy()
$ws.Range("M9").PasteSpecial($xlPasteFormats)
$ws.Range("M9").Value = 44.5

$ws.Range("K10").Copy()
$ws.Range("M10").PasteSpecial($xlPasteFormats)
$ws.Range("M10").Value = 21.9

$ws.Range("K11").Copy()
$ws.Range("M11").PasteSpecial($xlPasteFormats)
$ws.Range("M11").Value = 21.9

$ws.Range("K12").Copy()
$ws.Range("M12").PasteSpecial($xlPasteFormats)
$ws.Range("M12").Value = 35.3

$ws.Range("K13").Copy()
$ws.Range("M13").PasteSpecial($xlPasteFormats)
$ws.Range("M13").Value = 17.6

$ws.Range("K15").Copy()
$ws.Range("M15").PasteSpecial($xlPasteFormats)
$ws.Range("M15").Value = 44.9

$ws.Range("K16").Copy()
$ws.Range("M16").PasteSpecial($xlPasteFormats)
$ws.Range("M16").Value = 21.5

$ws.Range("K17").Copy()
$ws.Range("M17").PasteSpecial($xlPasteFormats)
$ws.Range("M17").Value = 36.2

$ws.Range("K18").Copy()
$ws.Range("M18").PasteSpecial($xlPasteFormats)
$ws.Range("M18").Value = 53.1

$ws.Range("K19").Copy()
$ws.Range("M19").PasteSpecial($xlPasteFormats)
$ws.Range("M19").Value = 33.4

$ws.Range("K20").Copy()
$ws.Range("M20").PasteSpecial($xlPasteFormats)
$ws.Range("M20").Value = 20.2

$ws.Range("K21").Copy()
$ws.Range("M21").PasteSpecial($xlPasteFormats)
$ws.Range("M21").Value = 15.4

$ws.Range("K22").Copy()
$ws.Range("M22").PasteSpecial($xlPasteFormats)
$ws.Range("M22").Value = 20.5

$ws.Range("K23").Copy()
$ws.Range("M23").PasteSpecial($xlPasteFormats)
$ws.Range("M23").Value = 27.1

$ws.Range("K24").Copy()
$ws.Range("M24").PasteSpecial($xlPasteFormats)
$ws.Range("M24").Value = 36.1

$ws.Range("K25").Copy()
$ws.Range("M25").PasteSpecial($xlPasteFormats)
$ws.Range("M25").Value = 25.2

$ws.Range("K26").Copy()
$ws.Range("M26").PasteSpecial($xlPasteFormats)
$ws.Range("M26").Value = 24.2

$ws.Range("K27").Copy()
$ws.Range("M27").PasteSpecial($xlPasteFormats)
$ws.Range("M27").Value = 46.5

$ws.Range("K28").Copy()
$ws.Range("M28").PasteSpecial($xlPasteFormats)
$ws.Range("M28").Value = 20.3

$ws.Range("K29").Copy()
$ws.Range("M29").PasteSpecial($xlPasteFormats)
$ws.Range("M29").Value = 40.5

# Row 14 needs an explicit numeric format (0.0) applied on top of the K14 style
$ws.Range("K14").Copy()
$ws.Range("M14").PasteSpecial($xlPasteFormats)
$ws.Range("M14").NumberFormat = "0.0"
$ws.Range("M14").Value = 28

# Row 30 (Bishkek city total) - same style as L30 (bordered + number format)
$ws.Range("L30").Copy()
$ws.Range("M30").PasteSpecial($xlPasteFormats)
$ws.Range("M30").Value = 44.5

$excel.CutCopyMode = $false

# Update the active selection on the sheet to match the final authored state
$ws.Range("N7").Select()
